$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting (avoid Excel
# auto-converting numeric-looking strings into numbers and losing trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.998.17"
$ws.Range("E2").Value = "  -0.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.762.47"
$ws.Range("E3").Value = "  -3.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.80"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("E7").Value = "  -4.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3353"
$ws.Range("E8").Value = "  -4.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.55"
$ws.Range("E9").Value = "  -5.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.119"
$ws.Range("E10").Value = "  -7.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07190"
$ws.Range("E11").Value = "  -5.62%  "

$ws.Range("E12").Value = "  -0.31%  "

$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.194"
$ws.Range("E14").Value = "  -5.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.183"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.756.82"
$ws.Range("E16").Value = "  -3.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001054"
$ws.Range("E17").Value = "  -5.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06580"
$ws.Range("E18").Value = "  -2.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.44"
$ws.Range("E19").Value = "  -6.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9995"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.98"
$ws.Range("E21").Value = "  -5.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.288"
$ws.Range("E22").Value = "  -5.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.990.28"
$ws.Range("E23").Value = "  -0.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.70"
$ws.Range("E24").Value = "  -7.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.355"
$ws.Range("E25").Value = "  -2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.39"
$ws.Range("E26").Value = "  -1.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.83"
$ws.Range("E27").Value = "  -7.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.339"
$ws.Range("E28").Value = "  -9.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.958.78"
$ws.Range("E29").Value = "  -3.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("E32").Value = "  -0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.808"
$ws.Range("E33").Value = "  -6.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08773"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.26"
$ws.Range("E35").Value = "  -7.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02339"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6591"
$ws.Range("E37").Value = "  -6.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06184"
$ws.Range("E38").Value = "  -6.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.151"
$ws.Range("E39").Value = "  -7.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2110"
$ws.Range("E40").Value = "  -5.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.211"
$ws.Range("E41").Value = "  -4.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.447"
$ws.Range("E42").Value = "  -10.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.008"
$ws.Range("E43").Value = "  -6.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9987"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.81"
$ws.Range("E45").Value = "  -6.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.822"
$ws.Range("E46").Value = "  -1.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6051"
$ws.Range("E47").Value = "  -7.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.20"
$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.012"
$ws.Range("E49").Value = "  -7.78%  "

$ws.Range("E50").Value = "  +1.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07161"
$ws.Range("E51").Value = "  -0.98%  "

# Swap rows 30 and 31 (ImmutableX / BitcoinCash reordering)
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.82"
$ws.Range("E30").Value = "  -3.52%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.257"
$ws.Range("E31").Value = "  -16.16%  "